$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column R: "uds. Objetivo semana pasada" bumped up for several rows
$ws.Range("R4").Value = 3
$ws.Range("R8").Value = 2
$ws.Range("R9").Value = 1
$ws.Range("R10").Value = 1
$ws.Range("R12").Value = 3
$ws.Range("R13").Value = 1
$ws.Range("R14").Value = 1
$ws.Range("R15").Value = 5
$ws.Range("R18").Value = 2
$ws.Range("R21").Value = 2
$ws.Range("R22").Value = 1
$ws.Range("R24").Value = 1
$ws.Range("R31").Value = 1

# Column T: "Tendencia Consumo" recalculated (S - R) for the rows above that have R>0 now
$ws.Range("T10").Value = 1
$ws.Range("T14").Value = 1
$ws.Range("T24").Value = 0
$ws.Range("T31").Value = 0

# Column U: "Pedido Final" recalculated for row 31
$ws.Range("U31").Value = 0

# Row 31 (3501010006 - VITHAL INSECTICIDA ACCION TOTAL 250ML) is now hidden
$ws.Rows(31).Hidden = $true

# Total_Unidades footer reflects the updated order total
$ws.Range("C38").Value = 41
